$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.269.03'
$ws.Range('E2').Value = '  -1.15%  '
$ws.Range('D3').Value = '2.241.92'
$ws.Range('E3').Value = '  -1.30%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '247.05'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.89%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.629'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.01%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '74.77'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -2.74%  '
$ws.Range('E8').Value = '  +0.11%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.618'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -3.20%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '42.27'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +6.31%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0945'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.31%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.09'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.59%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.104'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.96%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '14.52'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.86%  '
$ws.Range('E15').Value = '  -1.59%  '
$ws.Range('D16').Value = '2.243.27'
$ws.Range('E16').Value = '  -1.32%  '
$ws.Range('D17').Value = '42.087.89'
$ws.Range('E17').Value = '  -1.26%  '
$ws.Range('E18').Value = '  -0.28%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.16'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.16%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '72.18'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.23%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '231.94'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.36%  '
$ws.Range('E22').Value = '  +3.65%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.98'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +40.66%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.999'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.05%  '
$ws.Range('E25').Value = '  +0.24%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.63'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -4.38%  '
$ws.Range('E27').Value = '  -2.73%  '
$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.16'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.67%  '
$ws.Range('B29').Value = 'Monero'
$ws.Range('C29').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '168.99'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.96%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '20.66'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.93%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0821'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.57%  '
$ws.Range('B32').Value = 'InjectiveProtocol'
$ws.Range('C32').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '30.92'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.02%  '
$ws.Range('B33').Value = 'Kaspa'
$ws.Range('C33').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.120'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.67%  '
$ws.Range('E34').Value = '  -1.66%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.24'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +11.55%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.49'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.18%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0310'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.92%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '13.62'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.44%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.19'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.05%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.78'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.07%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '62.50'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.83%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.205'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.23%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '106.60'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.18%  '
$ws.Range('E44').Value = '  +2.39%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.69'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.60%  '
$ws.Range('E46').Value = '  -0.41%  '
$ws.Range('E47').Value = '  -2.65%  '
$ws.Range('B48').Value = 'TrustWalletToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.16'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.01%  '
$ws.Range('B49').Value = 'FTXToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '4.26'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -7.78%  '
$ws.Range('E50').Value = '  +1.79%  '
$ws.Range('E51').Value = '  +16.63%  '
